# Update the ACE AVIATION course row on the "courses" sheet.
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("courses")

# department (C2): "ACE AVIATION" -> "AVIATION"
$ws.Range("C2").Value = "AVIATION"

# promotionValidity (R2): remove the expired promo text
$ws.Range("R2").ClearContents()

# Leave the selection on the cell that was last edited.
[void]$ws.Range("R2").Select()
